$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.15"
$ws.Range("E2").Value = "'3.95%"
$ws.Range("D3").Value = "'32.04"
$ws.Range("E3").Value = "'8.48%"
$ws.Range("D4").Value = "'5.245"
$ws.Range("E4").Value = "'0.77%"
$ws.Range("D5").Value = "'0.07561"
$ws.Range("E5").Value = "'6.12%"
$ws.Range("D6").Value = "'7.844"
$ws.Range("E6").Value = "'4.40%"
$ws.Range("E7").Value = "'6.95%"
$ws.Range("D8").Value = "'1.536"
$ws.Range("E8").Value = "'9.12%"
$ws.Range("D9").Value = "'0.9280"
$ws.Range("E9").Value = "'1.87%"
$ws.Range("D10").Value = "'0.1692"
$ws.Range("E10").Value = "'4.02%"
$ws.Range("D11").Value = "'0.07938"
$ws.Range("E11").Value = "'3.55%"
$ws.Range("D12").Value = "'0.08065"
$ws.Range("E12").Value = "'4.46%"
$ws.Range("D13").Value = "'0.03018"
$ws.Range("E13").Value = "'2.41%"
$ws.Range("D14").Value = "'0.09913"
$ws.Range("E14").Value = "'10.01%"
$ws.Range("D15").Value = "'0.001490"
$ws.Range("E15").Value = "'-6.82%"
$ws.Range("D16").Value = "'0.04601"
$ws.Range("E16").Value = "'1.67%"
$ws.Range("D17").Value = "'0.006281"
$ws.Range("E17").Value = "'2.02%"
$ws.Range("E18").Value = "'-0.93%"
$ws.Range("D19").Value = "'2.233"
$ws.Range("E19").Value = "'-0.04%"
$ws.Range("D20").Value = "'0.3300"
$ws.Range("E20").Value = "'0.89%"
$ws.Range("D21").Value = "'0.1334"
$ws.Range("E21").Value = "'-2.37%"
$ws.Range("D22").Value = "'4.437"
$ws.Range("E22").Value = "'9.49%"
$ws.Range("E23").Value = "'1.14%"
$ws.Range("D24").Value = "'0.001214"
$ws.Range("E24").Value = "'0.30%"
$ws.Range("D25").Value = "'0.004474"
$ws.Range("E25").Value = "'5.15%"
$ws.Range("D26").Value = "'0.0001395"
$ws.Range("E26").Value = "'19.40%"
$ws.Range("D27").Value = "'0.0001783"
$ws.Range("E27").Value = "'5.55%"
$ws.Range("D39").Value = "'0.01712"
$ws.Range("E39").Value = "'2,508.75%"
$ws.Range("D40").Value = "'0.04500"
$ws.Range("E40").Value = "'2.51%"
$ws.Range("D41").Value = "'0.006960"
$ws.Range("E41").Value = "'-0.78%"
$ws.Range("D42").Value = "'0.1359"
$ws.Range("E42").Value = "'6.23%"
$ws.Range("D43").Value = "'0.002073"
$ws.Range("E43").Value = "'-6.08%"
$ws.Range("D44").Value = "'0.01376"
$ws.Range("E44").Value = "'1.81%"
$ws.Range("D45").Value = "'0.00006156"
$ws.Range("E45").Value = "'5.34%"
$ws.Range("D46").Value = "'1.845"
$ws.Range("E46").Value = "'-4.39%"
$ws.Range("D47").Value = "'0.01298"
$ws.Range("E47").Value = "'-0.09%"
